$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 799
$ws.Range("D4").Value = 0

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 27

$ws.Range("C6").Value = 788
$ws.Range("D6").Value = 5

$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 30

$ws.Range("C10").Value = 809
$ws.Range("D10").Value = 2

$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 16

$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 5

$ws.Range("C14").Value = 790
$ws.Range("D14").Value = 0

$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 32
